$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "2024-07-29 Monday" "2024-07-30 Tuesday"

Replace-Text "39×79=" "32×34="
Replace-Text "14×38=" "70×98="
Replace-Text "59×35=" "24×19="
Replace-Text "79×91=" "14×87="
Replace-Text "60×92=" "78×54="
Replace-Text "56×98=" "86×86="
Replace-Text "74×43=" "11×44="
Replace-Text "78×93=" "70×71="
Replace-Text "57×58=" "67×79="
Replace-Text "74×97=" "30×69="
Replace-Text "66×69=" "28×89="
Replace-Text "57×41=" "29×48="
Replace-Text "80×75=" "99×84="
Replace-Text "13×44=" "81×71="
Replace-Text "82×41=" "47×18="
Replace-Text "88×46=" "58×84="
Replace-Text "43×87=" "17×78="
Replace-Text "69×87=" "21×54="
Replace-Text "55×94=" "91×74="
Replace-Text "33×30=" "42×34="
Replace-Text "51×81=" "66×98="
Replace-Text "84×11=" "48×64="
Replace-Text "17×45=" "15×62="
Replace-Text "37×17=" "65×70="
Replace-Text "42×70=" "64×64="

$d.Save()
